# Actualizacion de horarios Linea 141 - 52
# Inserta nuevas filas de arribos en las hojas LP1912, LP1912-215 y 6203-6173,
# y actualiza los encabezados de "Ultima actualizacion" y "Total filas".
#
# Las filas nuevas se insertan de abajo hacia arriba (usando los numeros de
# fila originales como punto de anclaje) para que cada insercion no altere
# la posicion de los puntos de anclaje que todavia faltan procesar.

$wb = $excel.ActiveWorkbook

$nuevaHora = "Última actualización: 03:46:12"

function Insert-DataRow {
    param(
        $ws,
        [int]$rowNum,
        [string]$colA,
        [string]$colB,
        [string]$colC,
        [double]$colD
    )
    $ws.Rows.Item($rowNum).Insert()
    $ws.Cells.Item($rowNum, 1).Value = $colA
    $ws.Cells.Item($rowNum, 2).Value = $colB
    $ws.Cells.Item($rowNum, 3).Value = $colC
    $ws.Cells.Item($rowNum, 4).Value = $colD
}

# ---------------------------------------------------------------------------
# Hoja "LP1912" (y de forma identica "6203-6173"): 8 -> 12 filas de datos
# ---------------------------------------------------------------------------
$sheetNames = @("LP1912", "6203-6173")
foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Cells.Item(2, 1).Value = $nuevaHora
    $ws.Cells.Item(3, 1).Value = "Total filas: 12"

    # Ancla original fila 11 (81_EL PELIGRO): dos filas nuevas antes de ella.
    Insert-DataRow $ws 11 "03:46:12" "05:35" "215B_EL PATO" 109
    Insert-DataRow $ws 11 "03:46:12" "05:22" "23_HERNANDEZ" 96

    # Ancla original fila 9 (00:46:06 / 14_ABASTO): una fila nueva antes.
    Insert-DataRow $ws 9  "03:46:12" "05:16" "17_ROMERO" 90

    # Ancla original fila 7 (01:55:38 / 15_ABASTO): una fila nueva antes.
    Insert-DataRow $ws 7  "03:46:12" "04:46" "215A_EL PATO" 60
}

# ---------------------------------------------------------------------------
# Hoja "LP1912-215": 3 -> 5 filas de datos
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Cells.Item(2, 1).Value = $nuevaHora
$ws2.Cells.Item(3, 1).Value = "Total filas: 5"

# Ancla original fila 8 (215A_EL PATO / 118): una fila nueva antes.
Insert-DataRow $ws2 8 "03:46:12" "05:35" "215B_EL PATO" 109

# Ancla original fila 7 (215_ALUAR / 96): una fila nueva antes.
Insert-DataRow $ws2 7 "03:46:12" "04:46" "215A_EL PATO" 60
